$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 20, pushing existing rows 20-96 down to 22-98.
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# New row 20: Royal Gala
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = [DateTime]::new(2022, 2, 16)
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100104
$ws.Cells.Item(20, 8).Value = "Frutos de pepita"
$ws.Cells.Item(20, 9).Value = 100104002
$ws.Cells.Item(20, 10).Value = "Manzana"
$ws.Cells.Item(20, 11).Value = "Royal Gala"
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 300
$ws.Cells.Item(20, 14).Value = 21000
$ws.Cells.Item(20, 15).Value = 22000
$ws.Cells.Item(20, 16).Value = 21500
$ws.Cells.Item(20, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 1194
$ws.Cells.Item(20, 20).Value = 18

# New row 21: Scarlett
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = [DateTime]::new(2022, 2, 16)
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100104
$ws.Cells.Item(21, 8).Value = "Frutos de pepita"
$ws.Cells.Item(21, 9).Value = 100104002
$ws.Cells.Item(21, 10).Value = "Manzana"
$ws.Cells.Item(21, 11).Value = "Scarlett"
$ws.Cells.Item(21, 12).Value = "Segunda"
$ws.Cells.Item(21, 13).Value = 270
$ws.Cells.Item(21, 14).Value = 21000
$ws.Cells.Item(21, 15).Value = 22000
$ws.Cells.Item(21, 16).Value = 21500
$ws.Cells.Item(21, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(21, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 19).Value = 1194
$ws.Cells.Item(21, 20).Value = 18
